$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H129").Value = 3732840.2
$ws.Range("I129").Value = 35715316
$ws.Range("J129").Value = 1551.2
$ws.Range("K129").Value = 107145948
$ws.Range("L129").Value = 4653.6
$ws.Range("M129").Value = -107140948
$ws.Range("N129").Value = -14653.6

$ws.Range("H132").Value = 3775566.8
$ws.Range("I132").Value = 4168238.2
$ws.Range("J132").Value = 5919.8
$ws.Range("K132").Value = 12504714.6
$ws.Range("L132").Value = 17759.4
$ws.Range("M132").Value = -12502184.6
$ws.Range("N132").Value = -22819.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17859186
$ws.Range("I2").Value = 250000000
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 250000000
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -249999887

$ws.Range("H32").Value = 3477.716
$ws.Range("I32").Value = 2644.7246
$ws.Range("J32").Value = 8267.416999999999
$ws.Range("K32").Value = 2644.7246
$ws.Range("L32").Value = 8267.416999999999
$ws.Range("M32").Value = -2357.7246

$ws.Range("H61").Value = 3556.762
$ws.Range("I61").Value = 1441
$ws.Range("J61").Value = 6377.778
$ws.Range("K61").Value = 1441
$ws.Range("L61").Value = 6377.778
$ws.Range("M61").Value = -1229

$ws.Range("H74").Value = 913.1875
$ws.Range("I74").Value = 674.1818
$ws.Range("J74").Value = 1439
$ws.Range("K74").Value = 674.1818
$ws.Range("L74").Value = 1439
$ws.Range("M74").Value = 199.8182

$ws.Range("H77").Value = 913.1875
$ws.Range("I77").Value = 674.1818
$ws.Range("J77").Value = 1439
$ws.Range("K77").Value = 3370.909
$ws.Range("L77").Value = 7195
$ws.Range("M77").Value = 997.0910000000003

$ws.Range("H102").Value = 2742.6924
$ws.Range("I102").Value = 2404.7827
$ws.Range("J102").Value = 5333.3335
$ws.Range("K102").Value = 2404.7827
$ws.Range("L102").Value = 5333.3335
$ws.Range("M102").Value = -782.7827000000002

$ws.Range("H116").Value = 17859186
$ws.Range("I116").Value = 250000000
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 250000000
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -249997706

$ws.Range("H131").Value = 31666.666
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 31666.666
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 31666.666
$ws.Range("N131").Value = -41746.666

$ws.Range("H132").Value = 18870036
$ws.Range("I132").Value = 23811316
$ws.Range("J132").Value = 3335.4546
$ws.Range("K132").Value = 71433948
$ws.Range("L132").Value = 10006.3638
$ws.Range("M132").Value = -71431418
$ws.Range("N132").Value = -15066.3638

$ws.Range("H136").Value = 3556.762
$ws.Range("I136").Value = 1441
$ws.Range("J136").Value = 6377.778
$ws.Range("K136").Value = 4323
$ws.Range("L136").Value = 19133.334
$ws.Range("M136").Value = -1773

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17859186
$ws.Range("I3").Value = 250000000
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 250000000
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -249999886

$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372

$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864

$ws.Range("H82").Value = 13563.3
$ws.Range("I82").Value = 3235.7144
$ws.Range("J82").Value = 37661
$ws.Range("K82").Value = 3235.7144
$ws.Range("L82").Value = 37661
$ws.Range("M82").Value = -2852.7144
$ws.Range("N82").Value = -38427

$ws.Range("H85").Value = 13563.3
$ws.Range("I85").Value = 3235.7144
$ws.Range("J85").Value = 37661
$ws.Range("K85").Value = 3235.7144
$ws.Range("L85").Value = 37661
$ws.Range("M85").Value = -1909.7144
$ws.Range("N85").Value = -40313

$ws.Range("H134").Value = 3459.6365
$ws.Range("I134").Value = 2247
$ws.Range("J134").Value = 7582.6
$ws.Range("K134").Value = 6741
$ws.Range("L134").Value = 22747.8
$ws.Range("M134").Value = -4206
$ws.Range("N134").Value = -27817.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4020800.5
$ws.Range("I3").Value = 11666.667
$ws.Range("J3").Value = 10034502
$ws.Range("K3").Value = 11666.667
$ws.Range("L3").Value = 10034502
$ws.Range("M3").Value = -11553.667
$ws.Range("N3").Value = -10034728

$ws.Range("H7").Value = 94
$ws.Range("I7").Value = 55
$ws.Range("J7").Value = 152.5
$ws.Range("K7").Value = 55
$ws.Range("L7").Value = 152.5
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = -378.5

$ws.Range("H31").Value = 2098.2834
$ws.Range("I31").Value = 1270.4324
$ws.Range("J31").Value = 3430.0435
$ws.Range("K31").Value = 1270.4324
$ws.Range("L31").Value = 3430.0435
$ws.Range("M31").Value = -975.4323999999999
$ws.Range("N31").Value = -4020.0435

$ws.Range("H34").Value = 2098.2834
$ws.Range("I34").Value = 1270.4324
$ws.Range("J34").Value = 3430.0435
$ws.Range("K34").Value = 1270.4324
$ws.Range("L34").Value = 3430.0435
$ws.Range("M34").Value = -1068.4324
$ws.Range("N34").Value = -3834.0435

$ws.Range("H58").Value = 16669187
$ws.Range("I58").Value = 1443.8889
$ws.Range("J58").Value = 41670800
$ws.Range("K58").Value = 1443.8889
$ws.Range("L58").Value = 41670800
$ws.Range("M58").Value = -1240.8889
$ws.Range("N58").Value = -41671206

$ws.Range("H107").Value = 1773.16
$ws.Range("I107").Value = 426.66666
$ws.Range("J107").Value = 3016.077
$ws.Range("K107").Value = 426.66666
$ws.Range("L107").Value = 3016.077
$ws.Range("M107").Value = 1493.33334
$ws.Range("N107").Value = -6856.077

$ws.Range("H132").Value = 2988.0715
$ws.Range("I132").Value = 2142.5557
$ws.Range("J132").Value = 4510
$ws.Range("K132").Value = 6427.6671
$ws.Range("L132").Value = 13530
$ws.Range("M132").Value = -3897.6671
$ws.Range("N132").Value = -18590

$ws.Range("H134").Value = 2108.7144
$ws.Range("I134").Value = 1310.0435
$ws.Range("J134").Value = 5782.6
$ws.Range("K134").Value = 3930.1305
$ws.Range("L134").Value = 17347.8
$ws.Range("M134").Value = -1395.1305
$ws.Range("N134").Value = -22417.8

$ws.Range("H136").Value = 16669187
$ws.Range("I136").Value = 1443.8889
$ws.Range("J136").Value = 41670800
$ws.Range("K136").Value = 4331.6667
$ws.Range("L136").Value = 125012400
$ws.Range("M136").Value = -1781.6667
$ws.Range("N136").Value = -125017500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 955.7963
$ws.Range("I107").Value = 583.0476
$ws.Range("J107").Value = 1193
$ws.Range("K107").Value = 1749.1428
$ws.Range("L107").Value = 3579
$ws.Range("M107").Value = 170.8571999999999
$ws.Range("N107").Value = -7419

$ws.Range("H129").Value = 53593
$ws.Range("I129").Value = 5390
$ws.Range("J129").Value = 101796
$ws.Range("K129").Value = 16170
$ws.Range("L129").Value = 305388
$ws.Range("M129").Value = -11170
$ws.Range("N129").Value = -315388

$ws.Range("H131").Value = 1398.9697
$ws.Range("I131").Value = 883.3
$ws.Range("J131").Value = 2192.3076
$ws.Range("K131").Value = 2649.9
$ws.Range("L131").Value = 6576.9228
$ws.Range("M131").Value = 2390.1
$ws.Range("N131").Value = -16656.9228

$ws.Range("H137").Value = 3793.9583
$ws.Range("I137").Value = 3729.923
$ws.Range("J137").Value = 3869.6365
$ws.Range("K137").Value = 11189.769
$ws.Range("L137").Value = 11608.9095
$ws.Range("M137").Value = -6089.769
$ws.Range("N137").Value = -21808.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 49999
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 49999
$ws.Range("N42").Value = -50969

$ws.Range("H115").Value = 49999
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 49999
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 49999
$ws.Range("N115").Value = -52349

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 38464890
$ws.Range("I132").Value = 62502092
$ws.Range("J132").Value = 5359
$ws.Range("K132").Value = 187506276
$ws.Range("L132").Value = 16077
$ws.Range("M132").Value = -187503746
$ws.Range("N132").Value = -21137

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4285.2856
$ws.Range("I40").Value = 1555.7142
$ws.Range("J40").Value = 7014.857
$ws.Range("K40").Value = 1555.7142
$ws.Range("L40").Value = 7014.857
$ws.Range("M40").Value = -1419.7142
$ws.Range("N40").Value = -7286.857

$ws.Range("H122").Value = 3020.1333
$ws.Range("I122").Value = 2520
$ws.Range("J122").Value = 4020.4
$ws.Range("K122").Value = 7560
$ws.Range("L122").Value = 12061.2
$ws.Range("M122").Value = -5110
$ws.Range("N122").Value = -16961.2

$ws.Range("H136").Value = 1953.0322
$ws.Range("I136").Value = 1549.579
$ws.Range("J136").Value = 2591.8333
$ws.Range("K136").Value = 4648.737
$ws.Range("L136").Value = 7775.499899999999
$ws.Range("M136").Value = -2098.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2167.3684
$ws.Range("I122").Value = 1687.0588
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 5061.1764
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -2611.1764
$ws.Range("N122").Value = -23650

$ws.Range("H126").Value = 2453.7827
$ws.Range("I126").Value = 1888.3572
$ws.Range("J126").Value = 3333.3333
$ws.Range("K126").Value = 5665.071599999999
$ws.Range("L126").Value = 9999.999899999999
$ws.Range("M126").Value = -3195.071599999999
$ws.Range("N126").Value = -14939.9999

$ws.Range("H132").Value = 16887.428
$ws.Range("I132").Value = 2285.2173
$ws.Range("J132").Value = 44875
$ws.Range("K132").Value = 6855.651899999999
$ws.Range("L132").Value = 134625
$ws.Range("M132").Value = -4325.651899999999
$ws.Range("N132").Value = -139685
